# Auto-generated edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 2 (shifts existing rows 2-3 down to 4-5)
$ws.Rows("2:3").Insert()

# ---- Row 2 ----
$ws.Range("A2").Value = 1032170367
$ws.Range("B2").Value = 'Abbey'
$ws.Range("C2").Value = 'Paul'
$ws.Range("D2").Value = '-'
$ws.Range("E2").Value = '-'
$ws.Range("F2").Value = '-'
$ws.Range("G2").Value = 'Dentist'
$ws.Range("H2").Value = 1891768495
$ws.Range("I2").Value = '-'
$ws.Range("J2").Value = '-'
$ws.Range("K2").Value = '-'
$ws.Range("L2").Value = '-'
$ws.Range("M2").Value = '-'
$ws.Range("N2").Value = '-'
$ws.Range("O2").Value = '-'
$ws.Range("P2").Value = '-'
$ws.Range("Q2").Value = '-'
$ws.Range("R2").Value = '-'
$ws.Range("S2").Value = '-'
$ws.Range("T2").Value = '-'
$ws.Range("U2").Value = 'DEN.00104931'
$ws.Range("V2").Value = 'colorado'
$ws.Range("W2").Value = 'Paul Frederick Abbey'
$ws.Range("X2").Value = 'Active'
$ws.Range("Y2").Value = "'02/28/2022"
$ws.Range("Y2").ClearFormats()
$ws.Range("Z2").Value = 'No'
$ws.Range("AA2").Value = "'07/08/2020"
$ws.Range("AA2").ClearFormats()
$ws.Range("AB2").Value = '-'
$ws.Range("AC2").Value = '-'
$ws.Range("AD2").Value = '-'
$ws.Range("AE2").Value = '-'
$ws.Range("AF2").Value = '-'
$ws.Range("AG2").Value = '-'
$ws.Range("AH2").Value = '-'
$ws.Range("AI2").Value = '-'
$ws.Range("AJ2").Value = '-'
$ws.Range("AK2").Value = '-'
$ws.Range("AL2").Value = '-'
$ws.Range("AM2").Value = '-'
$ws.Range("AN2").Value = '-'
$ws.Range("AO2").Value = '-'
$ws.Range("AP2").Value = '-'
$ws.Range("AQ2").Value = '-'
$ws.Range("AR2").Value = '-'
$ws.Range("AS2").Value = '-'
$ws.Range("AT2").Value = '-'
$ws.Range("AU2").Value = '-'
$ws.Range("AV2").Value = '-'
$ws.Range("AW2").Value = '-'
$ws.Range("AX2").Value = '-'
$ws.Range("AY2").Value = '-'
$ws.Range("AZ2").Value = '-'
$ws.Range("BA2").Value = '-'
$ws.Range("BB2").Value = '-'
$ws.Range("BC2").Value = '-'
$ws.Range("BD2").Value = '-'
$ws.Range("BE2").Value = '-'
$ws.Range("BF2").Value = '-'
$ws.Range("BG2").Value = '-'
$ws.Range("BH2").Value = '-'
$ws.Range("BI2").Value = '-'
$ws.Range("BJ2").Value = '-'
$ws.Range("BK2").Value = '-'
$ws.Range("BL2").Value = '-'
$ws.Range("BM2").Value = '-'
$ws.Range("BN2").Value = '-'
$ws.Range("BO2").Value = '-'
$ws.Range("BP2").Value = '-'
$ws.Range("BQ2").Value = '-'
$ws.Range("BR2").Value = '-'
$ws.Range("BS2").Value = '-'
$ws.Range("BT2").Value = '-'
$ws.Range("BU2").Value = '-'
$ws.Range("BV2").Value = '-'
$ws.Range("BW2").Value = '-'
$ws.Range("BX2").Value = '-'
$ws.Range("BY2").Value = '-'
$ws.Range("BZ2").Value = '-'
$ws.Range("CA2").Value = '-'
$ws.Range("CB2").Value = '-'
$ws.Range("CC2").Value = '-'
$ws.Range("CD2").Value = '-'
$ws.Range("CE2").Value = '-'
$ws.Range("CF2").Value = '-'
$ws.Range("CG2").Value = '-'
$ws.Range("CH2").Value = '-'
$ws.Range("CI2").Value = '-'
$ws.Range("CJ2").Value = '-'
$ws.Range("CK2").Value = '-'
$ws.Range("CL2").Value = '-'
$ws.Range("CM2").Value = '-'
$ws.Range("CN2").Value = '-'
$ws.Range("CO2").Value = '-'
$ws.Range("CP2").Value = '-'
$ws.Range("CQ2").Value = '-'
$ws.Range("CR2").Value = '-'
$ws.Range("CS2").Value = '-'
$ws.Range("CT2").Value = '-'
$ws.Range("CU2").Value = '-'
$ws.Range("CV2").Value = '-'
$ws.Range("CW2").Value = '-'
$ws.Range("CX2").Value = '-'
$ws.Range("CY2").Value = '-'
$ws.Range("CZ2").Value = '-'
$ws.Range("DA2").Value = '-'
$ws.Range("DB2").Value = '-'
$ws.Range("DC2").Value = '-'
$ws.Range("DD2").Value = '-'
$ws.Range("DE2").Value = '-'
$ws.Range("DF2").Value = '-'
$ws.Range("DG2").Value = '-'
$ws.Range("DH2").Value = '-'
$ws.Range("DI2").Value = '-'
$ws.Range("DJ2").Value = '-'
$ws.Range("DK2").Value = '-'
$ws.Range("DL2").Value = '-'
$ws.Range("DM2").Value = '-'
$ws.Range("DN2").Value = '-'
$ws.Range("DO2").Value = '-'
$ws.Range("DP2").Value = '-'
$ws.Range("DQ2").Value = '-'
$ws.Range("DR2").Value = 'Paul F. Abbey'
$ws.Range("DS2").Value = 44075
$ws.Range("DT2").Value = '-'
$ws.Range("DU2").Value = 'no'
$ws.Range("DV2").Value = 'yes'
$ws.Range("DW2").Value = 'yes'
$ws.Range("DX2").Value = 'yes'
$ws.Range("DY2").Value = 'no'
$ws.Range("DZ2").Value = 'yes'
$ws.Range("EA2").Value = 'yes'
$ws.Range("EB2").Value = 'yes'
$ws.Range("EC2").Value = 'yes'
$ws.Range("ED2").Value = 'yes'
$ws.Range("EE2").Value = '-'
$ws.Range("EF2").Value = '-'
$ws.Range("EG2").Value = '-'
$ws.Range("EH2").Value = '-'
$ws.Range("EI2").Value = '-'
$ws.Range("EJ2").Value = '-'
$ws.Range("EK2").Value = '-'
$ws.Range("EL2").Value = '-'
$ws.Range("EM2").Value = '1507 W MOUNTAIN VIEW AVE'
$ws.Range("EN2").Value = 'LONGMONT'
$ws.Range("EO2").Value = 'CO'
$ws.Range("EP2").Value = '80501-3201'
$ws.Range("EQ2").Value = '-'
$ws.Range("ER2").Value = '-'
$ws.Range("ES2").Value = '-'
$ws.Range("ET2").Value = '-'
$ws.Range("EU2").Value = '-'
$ws.Range("EV2").Value = '-'
$ws.Range("EW2").Value = '-'
$ws.Range("EX2").Value = '-'
$ws.Range("EY2").Value = '-'
$ws.Range("EZ2").Value = '-'
$ws.Range("FA2").Value = '-'
$ws.Range("FB2").Value = '-'
$ws.Range("FC2").Value = '-'
$ws.Range("FD2").Value = '-'
$ws.Range("FE2").Value = '-'
$ws.Range("FF2").Value = '-'
$ws.Range("FG2").Value = '(303) 678-0997'
$ws.Range("FH2").Value = 'paul.abbey@gmail.com'
$ws.Range("FI2").Value = '-'

# ---- Row 3 ----
$ws.Range("A3").Value = 3432170234
$ws.Range("B3").Value = 'Abrams'
$ws.Range("C3").Value = 'Edward'
$ws.Range("D3").Value = '-'
$ws.Range("E3").Value = 34350
$ws.Range("F3").Value = '-'
$ws.Range("G3").Value = 'Dentist'
$ws.Range("H3").Value = 1568575868
$ws.Range("I3").Value = '-'
$ws.Range("J3").Value = '-'
$ws.Range("K3").Value = '-'
$ws.Range("L3").Value = '-'
$ws.Range("M3").Value = '-'
$ws.Range("N3").Value = '-'
$ws.Range("O3").Value = '-'
$ws.Range("P3").Value = '-'
$ws.Range("Q3").Value = '-'
$ws.Range("R3").Value = '-'
$ws.Range("S3").Value = '-'
$ws.Range("T3").Value = '-'
$ws.Range("U3").Value = 'DS025857L'
$ws.Range("V3").Value = 'pennstate'
$ws.Range("W3").Value = 'EDWARDABRAMS'
$ws.Range("X3").Value = "'3/31/2021"
$ws.Range("X3").ClearFormats()
$ws.Range("Y3").Value = 'Active'
$ws.Range("Z3").Value = 'No'
$ws.Range("AA3").Value = "'07/08/2020"
$ws.Range("AA3").ClearFormats()
$ws.Range("AB3").Value = '-'
$ws.Range("AC3").Value = '-'
$ws.Range("AD3").Value = '-'
$ws.Range("AE3").Value = '-'
$ws.Range("AF3").Value = '-'
$ws.Range("AG3").Value = '-'
$ws.Range("AH3").Value = '-'
$ws.Range("AI3").Value = '-'
$ws.Range("AJ3").Value = '-'
$ws.Range("AK3").Value = '-'
$ws.Range("AL3").Value = '-'
$ws.Range("AM3").Value = '-'
$ws.Range("AN3").Value = '-'
$ws.Range("AO3").Value = '-'
$ws.Range("AP3").Value = '-'
$ws.Range("AQ3").Value = '-'
$ws.Range("AR3").Value = '-'
$ws.Range("AS3").Value = '-'
$ws.Range("AT3").Value = '-'
$ws.Range("AU3").Value = '-'
$ws.Range("AV3").Value = '-'
$ws.Range("AW3").Value = '-'
$ws.Range("AX3").Value = '-'
$ws.Range("AY3").Value = '-'
$ws.Range("AZ3").Value = '-'
$ws.Range("BA3").Value = '-'
$ws.Range("BB3").Value = '-'
$ws.Range("BC3").Value = '-'
$ws.Range("BD3").Value = '-'
$ws.Range("BE3").Value = '-'
$ws.Range("BF3").Value = '-'
$ws.Range("BG3").Value = '-'
$ws.Range("BH3").Value = '-'
$ws.Range("BI3").Value = '-'
$ws.Range("BJ3").Value = '-'
$ws.Range("BK3").Value = '-'
$ws.Range("BL3").Value = '-'
$ws.Range("BM3").Value = '-'
$ws.Range("BN3").Value = '-'
$ws.Range("BO3").Value = '-'
$ws.Range("BP3").Value = '-'
$ws.Range("BQ3").Value = '-'
$ws.Range("BR3").Value = '-'
$ws.Range("BS3").Value = '-'
$ws.Range("BT3").Value = '-'
$ws.Range("BU3").Value = '-'
$ws.Range("BV3").Value = '-'
$ws.Range("BW3").Value = '-'
$ws.Range("BX3").Value = '-'
$ws.Range("BY3").Value = '-'
$ws.Range("BZ3").Value = '-'
$ws.Range("CA3").Value = '-'
$ws.Range("CB3").Value = '-'
$ws.Range("CC3").Value = '-'
$ws.Range("CD3").Value = '-'
$ws.Range("CE3").Value = '-'
$ws.Range("CF3").Value = '-'
$ws.Range("CG3").Value = '-'
$ws.Range("CH3").Value = '-'
$ws.Range("CI3").Value = '-'
$ws.Range("CJ3").Value = '-'
$ws.Range("CK3").Value = '-'
$ws.Range("CL3").Value = '-'
$ws.Range("CM3").Value = '-'
$ws.Range("CN3").Value = '-'
$ws.Range("CO3").Value = '-'
$ws.Range("CP3").Value = '-'
$ws.Range("CQ3").Value = '-'
$ws.Range("CR3").Value = '-'
$ws.Range("CS3").Value = '-'
$ws.Range("CT3").Value = '-'
$ws.Range("CU3").Value = '-'
$ws.Range("CV3").Value = '-'
$ws.Range("CW3").Value = '-'
$ws.Range("CX3").Value = '-'
$ws.Range("CY3").Value = '-'
$ws.Range("CZ3").Value = '-'
$ws.Range("DA3").Value = '-'
$ws.Range("DB3").Value = '-'
$ws.Range("DC3").Value = '-'
$ws.Range("DD3").Value = '-'
$ws.Range("DE3").Value = '-'
$ws.Range("DF3").Value = '-'
$ws.Range("DG3").Value = '-'
$ws.Range("DH3").Value = '-'
$ws.Range("DI3").Value = '-'
$ws.Range("DJ3").Value = '-'
$ws.Range("DK3").Value = '-'
$ws.Range("DL3").Value = '-'
$ws.Range("DM3").Value = '-'
$ws.Range("DN3").Value = '-'
$ws.Range("DO3").Value = '-'
$ws.Range("DP3").Value = '-'
$ws.Range("DQ3").Value = '-'
$ws.Range("DR3").Value = 'Edward S. Abrams'
$ws.Range("DS3").Value = 43986
$ws.Range("DT3").Value = '-'
$ws.Range("DU3").Value = 'yes'
$ws.Range("DV3").Value = 'yes'
$ws.Range("DW3").Value = 'yes'
$ws.Range("DX3").Value = 'yes'
$ws.Range("DY3").Value = 'yes'
$ws.Range("DZ3").Value = 'no'
$ws.Range("EA3").Value = 'no'
$ws.Range("EB3").Value = 'no'
$ws.Range("EC3").Value = 'no'
$ws.Range("ED3").Value = 'no'
$ws.Range("EE3").Value = '-'
$ws.Range("EF3").Value = '-'
$ws.Range("EG3").Value = '-'
$ws.Range("EH3").Value = '-'
$ws.Range("EI3").Value = '-'
$ws.Range("EJ3").Value = '-'
$ws.Range("EK3").Value = '-'
$ws.Range("EL3").Value = '-'
$ws.Range("EM3").Value = '2137 WELSH ROAD, STE 3A'
$ws.Range("EN3").Value = 'PHILADELPHIA'
$ws.Range("EO3").Value = 'PA'
$ws.Range("EP3").Value = 19115
$ws.Range("EQ3").Value = '-'
$ws.Range("ER3").Value = '-'
$ws.Range("ES3").Value = '-'
$ws.Range("ET3").Value = '-'
$ws.Range("EU3").Value = '-'
$ws.Range("EV3").Value = '-'
$ws.Range("EW3").Value = '-'
$ws.Range("EX3").Value = '-'
$ws.Range("EY3").Value = '-'
$ws.Range("EZ3").Value = '-'
$ws.Range("FA3").Value = '-'
$ws.Range("FB3").Value = '-'
$ws.Range("FC3").Value = '-'
$ws.Range("FD3").Value = '-'
$ws.Range("FE3").Value = '-'
$ws.Range("FF3").Value = '-'
$ws.Range("FG3").Value = '215-969-1222'
$ws.Range("FH3").Value = 'abrams.edward@gmail.com'
$ws.Range("FI3").Value = '-'

